$wb = $excel.ActiveWorkbook

# --- Add the new sheet by copying "Namrata 105" (same row/column layout family:
# single line item + Total/GST/Grand Total rows) to the very end of the workbook ---
$src = $wb.Worksheets.Item("Namrata 105")
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$src.Copy($null, $last)

$ws = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws.Name = "Namrata 111"

# --- Update the quantity for the single line item; dependent formulas recalc ---
$ws.Range("D2").Value = 20

# --- Match the column widths / row heights seen on the new sheet ---
# (ColumnWidth is in "characters"; the host snaps to a whole-pixel grid, so the
# inputs below are chosen to land on the exact stored widths from the sheet.)
$ws.Range("A1").ColumnWidth = 6.5
$ws.Range("B1").ColumnWidth = 12.333333333333334
$ws.Range("C1").ColumnWidth = 22
$ws.Range("E1").ColumnWidth = 15.166666666666666

$ws.Rows.Item(1).RowHeight = 23.4
$ws.Rows.Item(2).RowHeight = 48

# --- Selection / active sheet bookkeeping ---
$rashi = $wb.Worksheets.Item("Rashi 109")
$rashi.Activate()
$rashi.Range("I24").Select()

$ws.Activate()
